$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-08-14 Wednesday" "2024-08-15 Thursday"

Replace-Text "837×8=6696" "959×8=7672"
Replace-Text "419×8=3352" "234×5=1170"
Replace-Text "248×4=992" "301×2=602"
Replace-Text "676×8=5408" "809×8=6472"
Replace-Text "365×5=1825" "795×6=4770"

Replace-Text "207×2=414" "234×5=1170"
Replace-Text "115×4=460" "580×4=2320"
Replace-Text "389×5=1945" "280×2=560"
Replace-Text "780×9=7020" "260×9=2340"
Replace-Text "457×4=1828" "931×7=6517"

Replace-Text "128×2=256" "944×4=3776"
Replace-Text "492×2=984" "569×4=2276"
Replace-Text "616×4=2464" "684×8=5472"
Replace-Text "177×8=1416" "761×2=1522"
Replace-Text "376×5=1880" "406×5=2030"

Replace-Text "665×5=3325" "727×4=2908"
Replace-Text "210×6=1260" "186×6=1116"
Replace-Text "919×8=7352" "764×2=1528"
Replace-Text "570×8=4560" "687×9=6183"
Replace-Text "674×4=2696" "548×2=1096"

Replace-Text "841×7=5887" "870×4=3480"
Replace-Text "434×5=2170" "108×7=756"
Replace-Text "668×3=2004" "473×3=1419"
Replace-Text "820×6=4920" "877×8=7016"
Replace-Text "746×8=5968" "513×7=3591"
